$d = $word.ActiveDocument

# 1. Update the DATE field from 2024-07-16 to 2024-07-25
$d.Content.Find.Execute("2024-07-16", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-25", 2) | Out-Null

# Helper: insert a new BodyText-styled bibliography paragraph immediately
# before the paragraph whose text currently starts with $anchorText.
function Insert-RefBefore($anchorText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $target = $rng.Duplicate
    $target.Collapse(1)
    $target.Text = $newText + "`r"
}

# 2. BC Housing entry, before the BC Stats / BC Demographic Survey entry
Insert-RefBefore "BC Stats. [creator] (2023). BC Demographic Survey." "BC Housing. [creator] (2023). Private Market Rent Supplements. E05. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

# 3. Human Early Learning Partnership entry, before Child Care Subsidy entry
Insert-RefBefore "Ministry of Children and Family Development. [creator] (2022). Child Care Subsidy." "Human Early Learning Partnership. [creator] (2024). Early Development Instrument. E03. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

# 4. Affordable Child Care Benefit entry, before K to 12 Socio Economic Status Index entry
Insert-RefBefore "Ministry of Education and Child Care. [creator] (2023). K to 12 Socio Economic Status Index." "Ministry of Education and Child Care. [creator] (2023). Affordable Child Care Benefit. E02. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

# 5. Neighbourhood Income entry, before BC Vital Events and Statistics entry
Insert-RefBefore "Ministry of Health. [creator] (2022). BC Vital Events and Statistics." "Ministry of Finance. [creator] (2024). Neighbourhood Income. E01. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

# 6. Home and Community Care entry, before Hospital Discharges entry
Insert-RefBefore "Ministry of Health. [creator] (2023). Hospital Discharges." "Ministry of Health. [creator] (2024). Home and Community Care. E03. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

# 7. PharmaCare and PharmaNet entries, before Registration and Premium Billings entry
Insert-RefBefore "Ministry of Health. [creator] (2019). Registration and Premium Billings." "Ministry of Health. [creator] (2020). PharmaCare. E01. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
Insert-RefBefore "Ministry of Health. [creator] (2019). Registration and Premium Billings." "Ministry of Health. [creator] (2023). PharmaNet. E02. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."

Write-Output "edits applied"
